$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.999999992072653
$ws.Range("A2").Value = 0.99520843136657289
$ws.Range("A3").Value = 0.97610727898855565
$ws.Range("A4").Value = 0.96787899161395186
$ws.Range("A5").Value = 0.96011651837320333
$ws.Range("A6").Value = 0.94440419631916495
$ws.Range("A7").Value = 0.94129012382374588
$ws.Range("A8").Value = 0.93789450662011498
$ws.Range("A9").Value = 0.93577963263632791
$ws.Range("A10").Value = 0.93460903573631582
$ws.Range("A11").Value = 0.93457051277186953
$ws.Range("A12").Value = 0.93475366944524774
$ws.Range("A13").Value = 0.93732856909643436
$ws.Range("A14").Value = 0.93923038544225323
$ws.Range("A15").Value = 0.94240407985251107
$ws.Range("A16").Value = 0.94062387646810264
$ws.Range("A17").Value = 0.93932450290834124
$ws.Range("A18").Value = 0.93898662431469415
$ws.Range("A19").Value = 0.99330877442662335
$ws.Range("A20").Value = 0.97841296835723046
$ws.Range("A21").Value = 0.97701452052830418
$ws.Range("A22").Value = 0.97575003468332289
$ws.Range("A23").Value = 0.96907194053331724
$ws.Range("A24").Value = 0.95605050919933088
$ws.Range("A25").Value = 0.94959343118430461
$ws.Range("A26").Value = 0.94236924072273776
$ws.Range("A27").Value = 0.93978445594611171
$ws.Range("A28").Value = 0.92398977416679573
$ws.Range("A29").Value = 0.91290542548722875
$ws.Range("A30").Value = 0.90780572229901912
$ws.Range("A31").Value = 0.90586874971804876
$ws.Range("A32").Value = 0.90635628399940482
$ws.Range("A33").Value = 0.90583629366458929
